# Apply updates described by the diff:
#  - Tweak a few odds values in row 5 (L5, M5, O5, Q5)
#  - Append a brand new data row (row 7) with match/odds info,
#    extending the sheet dimension from A1:BD6 to A1:BD7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 odds corrections ---
$ws.Range("L5").Value = 2.87
$ws.Range("M5").Value = 1.03
$ws.Range("O5").Value = 1.19
$ws.Range("Q5").Value = 1.75

# --- New row 7 ---
# Text / identifier columns
$ws.Range("A7").Value = "0CEt6wPQ"
$ws.Range("B7").Value = "21/11/2024"
$ws.Range("C7").Value = "12:30"
$ws.Range("D7").Value = "UNITED ARAB EMIRATES - UAE LEAGUE"
$ws.Range("E7").Value = "Al Urooba"
$ws.Range("F7").Value = "Al Ain"

# Numeric odds columns (G7:BC7)
$numericValues = @{
    "G7" = 9
    "H7" = 4.75
    "I7" = 1.32
    "J7" = 7
    "K7" = 2.57
    "L7" = 1.72
    "M7" = 1.02
    "N7" = 10
    "O7" = 1.11
    "P7" = 5.5
    "Q7" = 1.37
    "R7" = 2.85
    "S7" = 1.23
    "T7" = 3.75
    "U7" = 1.57
    "V7" = 2.25
    "W7" = 40
    "X7" = 90
    "Y7" = 27
    "Z7" = 300
    "AA7" = 90
    "AB7" = 50
    "AC7" = 10
    "AD7" = 10.5
    "AE7" = 15
    "AF7" = 45
    "AG7" = 200
    "AH7" = 10.75
    "AI7" = 8.75
    "AJ7" = 8.25
    "AK7" = 10
    "AL7" = 9.75
    "AM7" = 17.5
    "AN7" = 10.25
    "AO7" = 45
    "AP7" = 32
    "AQ7" = 300
    "AR7" = 200
    "AS7" = 250
    "AT7" = 3.75
    "AU7" = 6.9
    "AV7" = 40
    "AW7" = 3.55
    "AX7" = 6
    "AY7" = 11.75
    "AZ7" = 14.5
    "BA7" = 29
    "BB7" = 100
    "BC7" = 500
}

foreach ($ref in $numericValues.Keys) {
    $ws.Range($ref).Value = $numericValues[$ref]
}

# BD7 stays an empty (inline string) cell, mirroring BD6.
# Assigning Value = "" does not materialize a cell, so copy the existing
# empty BD6 cell down to BD7 to create an empty cell entry there.
$ws.Range("BD6").Copy($ws.Range("BD7"))
$excel.CutCopyMode = $false
